$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.14757047140632
$ws.Range("C2").Value = -1.292459396593679
$ws.Range("D2").Value = -1.35363339659368
$ws.Range("E2").Value = 0.6878226034063205
$ws.Range("F2").Value = 0.1052276034063205
$ws.Range("G2").Value = 0.2009027474063205
$ws.Range("H2").Value = -0.05951639659367947
$ws.Range("I2").Value = 0.7759386034063205
$ws.Range("J2").Value = 0.6758656504063205
$ws.Range("K2").Value = 0.3220726034063205
$ws.Range("B3").Value = -3.440029868393072
$ws.Range("C3").Value = -3.501203868393072
$ws.Range("D3").Value = -1.459747868393072
$ws.Range("E3").Value = -2.042342868393072
$ws.Range("F3").Value = -1.946667724393072
$ws.Range("G3").Value = -2.207086868393072
$ws.Range("H3").Value = -1.371631868393072
$ws.Range("I3").Value = -1.471704821393072
$ws.Range("J3").Value = -1.825497868393072
$ws.Range("K3").Value = -1.677319868393072
$ws.Range("B4").Value = -0.06117388186974204
$ws.Range("C4").Value = 1.980282118130258
$ws.Range("D4").Value = 1.397687118130258
$ws.Range("E4").Value = 1.493362262130258
$ws.Range("F4").Value = 1.232943118130258
$ws.Range("G4").Value = 2.068398118130258
$ws.Range("H4").Value = 1.968325165130258
$ws.Range("I4").Value = 1.614532118130258
$ws.Range("J4").Value = 1.762710118130258
$ws.Range("K4").Value = 1.920156118130258
$ws.Range("B5").Value = 2.041456296544459
$ws.Range("C5").Value = 1.458861296544459
$ws.Range("D5").Value = 1.554536440544459
$ws.Range("E5").Value = 1.294117296544459
$ws.Range("F5").Value = 2.129572296544459
$ws.Range("G5").Value = 2.029499343544459
$ws.Range("H5").Value = 1.675706296544459
$ws.Range("I5").Value = 1.823884296544459
$ws.Range("J5").Value = 1.981330296544459
$ws.Range("K5").Value = 1.480667296544459
$ws.Range("B6").Value = -0.5825946406117954
$ws.Range("C6").Value = -0.4869194966117955
$ws.Range("D6").Value = -0.7473386406117954
$ws.Range("E6").Value = 0.08811635938820461
$ws.Range("F6").Value = -0.01195659361179541
$ws.Range("G6").Value = -0.3657496406117954
$ws.Range("H6").Value = -0.2175716406117954
$ws.Range("I6").Value = -0.06012564061179543
$ws.Range("J6").Value = -0.5607886406117955
$ws.Range("K6").Value = -0.2804276406117954
$ws.Range("B7").Value = 0.09567504042184558
$ws.Range("C7").Value = -0.1647441035781544
$ws.Range("D7").Value = 0.6707108964218456
$ws.Range("E7").Value = 0.5706379434218456
$ws.Range("F7").Value = 0.2168448964218456
$ws.Range("G7").Value = 0.3650228964218456
$ws.Range("H7").Value = 0.5224688964218456
$ws.Range("I7").Value = 0.02180589642184558
$ws.Range("J7").Value = 0.3021668964218456
$ws.Range("K7").Value = 0.1301918964218456
$ws.Range("B8").Value = -0.2604191443875122
$ws.Range("C8").Value = 0.5750358556124878
$ws.Range("D8").Value = 0.4749629026124878
$ws.Range("E8").Value = 0.1211698556124878
$ws.Range("F8").Value = 0.2693478556124878
$ws.Range("G8").Value = 0.4267938556124878
$ws.Range("H8").Value = -0.07386914438751221
$ws.Range("I8").Value = 0.2064918556124878
$ws.Range("J8").Value = 0.0345168556124878
$ws.Range("K8").Value = 0.3403798556124878
$ws.Range("B9").Value = 0.8354548926112106
$ws.Range("C9").Value = 0.7353819396112107
$ws.Range("D9").Value = 0.3815888926112106
$ws.Range("E9").Value = 0.5297668926112107
$ws.Range("F9").Value = 0.6872128926112107
$ws.Range("G9").Value = 0.1865498926112106
$ws.Range("H9").Value = 0.4669108926112106
$ws.Range("I9").Value = 0.2949358926112106
$ws.Range("J9").Value = 0.6007988926112107
$ws.Range("K9").Value = -0.01153110738878937
$ws.Range("B10").Value = -0.1000730565472806
$ws.Range("C10").Value = -0.4538661035472806
$ws.Range("D10").Value = -0.3056881035472806
$ws.Range("E10").Value = -0.1482421035472806
$ws.Range("F10").Value = -0.6489051035472806
$ws.Range("G10").Value = -0.3685441035472806
$ws.Range("H10").Value = -0.5405191035472806
$ws.Range("I10").Value = -0.2346561035472806
$ws.Range("J10").Value = -0.8469861035472805
$ws.Range("K10").Value = -0.1586151035472806
$ws.Range("B11").Value = -0.3537867436446591
$ws.Range("C11").Value = -0.2056087436446591
$ws.Range("D11").Value = -0.04816274364465911
$ws.Range("E11").Value = -0.5488257436446591
$ws.Range("F11").Value = -0.2684647436446591
$ws.Range("G11").Value = -0.4404397436446591
$ws.Range("H11").Value = -0.1345767436446591
$ws.Range("I11").Value = -0.746906743644659
$ws.Range("J11").Value = -0.05853574364465908
$ws.Range("B12").Value = 0.1481777624350372
$ws.Range("C12").Value = 0.3056237624350371
$ws.Range("D12").Value = -0.1950392375649629
$ws.Range("E12").Value = 0.08532176243503714
$ws.Range("F12").Value = -0.08665323756496286
$ws.Range("G12").Value = 0.2192097624350371
$ws.Range("H12").Value = -0.3931202375649628
$ws.Range("I12").Value = 0.2952507624350372
$ws.Range("B13").Value = 0.1574463720025918
$ws.Range("C13").Value = -0.3432166279974082
$ws.Range("D13").Value = -0.06285562799740819
$ws.Range("E13").Value = -0.2348306279974082
$ws.Range("F13").Value = 0.0710323720025918
$ws.Range("G13").Value = -0.5412976279974082
$ws.Range("H13").Value = 0.1470733720025918
$ws.Range("B14").Value = -0.5006626170015632
$ws.Range("C14").Value = -0.2203016170015632
$ws.Range("D14").Value = -0.3922766170015632
$ws.Range("E14").Value = -0.08641361700156319
$ws.Range("F14").Value = -0.6987436170015632
$ws.Range("G14").Value = -0.01037261700156317
$ws.Range("B15").Value = 0.2803578395245076
$ws.Range("C15").Value = 0.1083828395245076
$ws.Range("D15").Value = 0.4142458395245076
$ws.Range("E15").Value = -0.1980841604754924
$ws.Range("F15").Value = 0.4902868395245076
$ws.Range("B16").Value = -0.1719750410109616
$ws.Range("C16").Value = 0.1338879589890384
$ws.Range("D16").Value = -0.4784420410109615
$ws.Range("E16").Value = 0.2099289589890385
$ws.Range("B17").Value = 0.3058628168340501
$ws.Range("C17").Value = -0.3064671831659499
$ws.Range("D17").Value = 0.3819038168340501
$ws.Range("B18").Value = -0.6123297229122814
$ws.Range("C18").Value = 0.07604127708771863
$ws.Range("B19").Value = 0.6883712297750049

$ws.Range("K11").ClearContents()
$ws.Range("J12").ClearContents()
$ws.Range("I13").ClearContents()
$ws.Range("H14").ClearContents()
$ws.Range("G15").ClearContents()
$ws.Range("F16").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("B20").ClearContents()

Write-Output "done"